$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Equipos")

# "Modelo" column (F) originally held the shared string "Hostia" for several
# rows (F2, F11, F12, F14, F15, F20, F21, F22, F23).
#
# For rows 20, 21 and 23 the value becomes "gg"; for row 22 it becomes "g".
# The remaining rows that still say "Hostia" are renamed to "aa".

$replaced = $ws.Cells.Replace("Hostia", "aa")

$ws.Range("F20").Value = "gg"
$ws.Range("F21").Value = "gg"
$ws.Range("F22").Value = "g"
$ws.Range("F23").Value = "gg"
